$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "SCOTT, Michael",
    "BERTRAM, Nellie",
    "CALIFORNIA, Robert",
    "LEVINSON, Jan",
    "MILLER, Pete",
    "ANDERSON, Roy",
    "MINER, Charles",
    "BENNETT, Jo",
    "GREEN, Clark",
    "VICKERS, Deangelo",
    "VANCE, Bob"
)

$startFileNumber = 123473
$startRow = 19

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startFileNumber + $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

$ws.Range("D20").Select()
